# Swap the "sessions" and "participants" sheet contents/names (keeping tab
# order/position the same), and update two header labels to include units.
#
# Before: tab1 = "sessions" (A1:AC1, session columns), tab2 = "participants"
#         (A1:J1, participant columns).
# After:  tab1 = "participants" (A1:J1, participant columns incl. updated
#         "gestational age (weeks)" / "birth weight (grams)" headers),
#         tab2 = "sessions" (A1:AC1, session columns) -- i.e. the sheet
#         contents are swapped between the two (still-in-place) tabs.

$wb = $excel.ActiveWorkbook

$sessionsSheet     = $wb.Worksheets.Item(1)   # currently "sessions"
$participantsSheet = $wb.Worksheets.Item(2)   # currently "participants"

# --- target header rows -----------------------------------------------
$sessionHeaders = @(
    "name", "key", "date", "participantID", "top", "pilot", "exclusion",
    "setting", "country", "state", "language", "release", "condition",
    "group", "tasks", "transcode_options", "filepath", "file_1", "fname_1",
    "fposition_1", "fclassification_1", "clip_out_1", "clip_in_1", "file_2",
    "fname_2", "fposition_2", "fclassification_2", "clip_out_2", "clip_in_2"
)

$participantHeaders = @(
    "participantID", "birthdate", "gender", "race", "ethnicity", "language",
    "disability", "gestational age (weeks)", "pregnancy term",
    "birth weight (grams)"
)

# --- wipe existing contents + validations on both sheets ---------------
foreach ($ws in @($sessionsSheet, $participantsSheet)) {
    $ws.Cells.Validation.Delete()
    $ws.Cells.ClearContents()
}

# --- rename via a temporary name to avoid a collision -------------------
$sessionsSheet.Name = "__tmp_swap__"
$participantsSheet.Name = "sessions"
$sessionsSheet.Name = "participants"

# sessionsSheet variable now refers to the tab that will hold "participants"
# data; participantsSheet variable now refers to the tab holding "sessions"
# data. Rename local references for clarity.
$newParticipantsTab = $sessionsSheet       # tab 1, renamed to "participants"
$newSessionsTab     = $participantsSheet   # tab 2, renamed to "sessions"

# --- write participant headers into tab 1 --------------------------------
for ($i = 0; $i -lt $participantHeaders.Length; $i++) {
    $newParticipantsTab.Cells.Item(1, $i + 1).Value = $participantHeaders[$i]
}

$newParticipantsTab.Range("C2:C1001").Validation.Add(3, 1, 1, "Female,Male")
$newParticipantsTab.Range("D2:D1001").Validation.Add(3, 1, 1, "American Indian or Alaska Native,Asian,Native Hawaiian or Other Pacific Islander,Black or African American,White,More than one,Unknown or not reported")
$newParticipantsTab.Range("E2:E1001").Validation.Add(3, 1, 1, "Not Hispanic or Latino,Hispanic or Latino")
$newParticipantsTab.Range("I2:I1001").Validation.Add(3, 1, 1, "Full term,Preterm")

# --- write session headers into tab 2 ------------------------------------
for ($i = 0; $i -lt $sessionHeaders.Length; $i++) {
    $newSessionsTab.Cells.Item(1, $i + 1).Value = $sessionHeaders[$i]
}

$newSessionsTab.Range("G2:G1001").Validation.Add(3, 1, 1, "Did not meet inclusion criteria,Procedural/experimenter error,Withdrew/fussy/tired,Outlier")
$newSessionsTab.Range("AA2:AA1001").Validation.Add(3, 1, 1, "None,PRIVATE,SHARED,EXCERPTS,PUBLIC")
$newSessionsTab.Range("J2:J1001").Validation.Add(3, 1, 1, "AL,AK,AZ,AR,CA,CO,CT,DE,DC,FL,GA,HI,ID,IL,IN,IA,KS,KY,LA,ME,MT,NE,NV,NH,NJ,NM,NY,NC,ND,OH,OK,OR,MD,MA,MI,MN,MS,MO,PA,RI,SC,SD,TN,TX,UT,VT,VA,WA,WV,WI,WY")
$newSessionsTab.Range("U2:U1001").Validation.Add(3, 1, 1, "None,PRIVATE,SHARED,EXCERPTS,PUBLIC")
$newSessionsTab.Range("L2:L1001").Validation.Add(3, 1, 1, "None,PRIVATE,SHARED,EXCERPTS,PUBLIC")
$newSessionsTab.Range("H2:H1001").Validation.Add(3, 1, 1, "Lab,Home,Classroom,Outdoor,Clinic")
